$wb = $excel.ActiveWorkbook

# 1. Rename the first sheet: "Complex Test Case" -> "Complex Test Case - Past"
$wsPast = $wb.Worksheets.Item("Complex Test Case")
$wsPast.Name = "Complex Test Case - Past"

# 2. "Complex Test Case - Past": selection B15:B20 -> B33, and give it a page setup
#    (paper size 9 / portrait) like the other sheets in the workbook.
$wsPast.Activate()
[void]$wsPast.Range("B33").Select()
$wsPast.PageSetup.PaperSize = 9
$wsPast.PageSetup.Orientation = 1

# 3. "Complex Test Case - Updated": selection D13 -> E25
$wsUpdated = $wb.Worksheets.Item("Complex Test Case - Updated")
$wsUpdated.Activate()
[void]$wsUpdated.Range("E25").Select()

# 4. "Unit Test Case": scrolled view A139 -> A61, selection F194 -> E77.
#    Keep this sheet active/selected last so it remains the active tab
#    (matches the unchanged activeTab="2" / tabSelected="1" in the workbook).
$wsUnit = $wb.Worksheets.Item("Unit Test Case")
$wsUnit.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
[void]$wsUnit.Range("E77").Select()
